$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column J (2021) data, mirroring the existing 2015-2020 (D:I) columns ---

# Row 3: blank separator cell under the header row, inherits the thick-bottom
# border formatting from I3 (same row).
$ws.Range("I3").Copy() | Out-Null
$ws.Range("J3").PasteSpecial(-4122) | Out-Null

# Row 4: year header value "2021", same look as the other year headers.
$ws.Range("J4").Value = 2021
$ws.Range("I4").Copy() | Out-Null
$ws.Range("J4").PasteSpecial(-4122) | Out-Null

# Rows 5 and 6 carry an explicit row-level style (customFormat), which a
# freshly-written cell in that row picks up automatically - no extra
# formatting step required.
$ws.Range("J5").Value = 5356.3
$ws.Range("J6").Value = 9.5

# Rows 8, 9, 11, 12, 14-17: plain data cells that pick up the worksheet's
# default column style automatically.
$ws.Range("J8").Value = 7.9
$ws.Range("J9").Value = 10.5
$ws.Range("J11").Value = 9.6
$ws.Range("J12").Value = 9.4
$ws.Range("J14").Value = 14.8
$ws.Range("J15").Value = 9.1
$ws.Range("J16").Value = 9.5
$ws.Range("J17").Value = 5.9

# Row 27: blank closing cell of the table, inherits the bottom-border style
# from I27 (same row), mirroring row 3's treatment at the top of the table.
$ws.Range("I27").Copy() | Out-Null
$ws.Range("J27").PasteSpecial(-4122) | Out-Null

# Leave the selection where the editor's cursor ended up.
$ws.Range("L27").Select() | Out-Null
